$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "지도 시각화 도구 Folium 사용법을 파헤쳐보자!"
$ws.Range("E4").Value = "https://teddylee777.github.io/visualization/folium"

$ws.Range("D9").Value = "MBA AI/BigData 첫 학기 시험 문제 공개"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/mba-ai-bigdata-1st-term-exam-sample/#utm_source=rss&utm_medium=rss&utm_campaign=mba-ai-bigdata-1st-term-exam-sample"

$ws.Range("D12").Value = "‘혼자 공부하는 머신러닝+딥러닝’이 Yes24 올해의 책 후보에 선정되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/11/08/%ed%98%bc%ec%9e%90-%ea%b3%b5%eb%b6%80%ed%95%98%eb%8a%94-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d%eb%94%a5%eb%9f%ac%eb%8b%9d%ec%9d%b4-yes24-%ec%98%ac%ed%95%b4%ec%9d%98-%ec%b1%85-%ed%9b%84%eb%b3%b4/"

$ws.Range("D29").Value = "[만화] 인턴일기 58~65"
$ws.Range("E29").Value = "https://blog.promedius.ai/intern-life-9/"
